# Auto-generated: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across several sheets to refreshed market-board values.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 121.0
$ws.Range("I5").Value = 112.28571
$ws.Range("J5").Value = 141.33333
$ws.Range("K5").Value = 112.28571
$ws.Range("L5").Value = 141.33333
$ws.Range("M5").Value = 2.714290000000005
$ws.Range("N5").Value = -371.33333
$ws.Range("H9").Value = 156.8
$ws.Range("J9").Value = 300.0
$ws.Range("L9").Value = 300.0
$ws.Range("N9").Value = -638.0
$ws.Range("H17").Value = 1741.5555
$ws.Range("J17").Value = 1814.5883
$ws.Range("L17").Value = 5443.7649
$ws.Range("N17").Value = -5779.7649
$ws.Range("H55").Value = 762.8125
$ws.Range("I55").Value = 518.4286
$ws.Range("J55").Value = 952.8889
$ws.Range("K55").Value = 518.4286
$ws.Range("L55").Value = 952.8889
$ws.Range("M55").Value = -304.4286
$ws.Range("N55").Value = -1380.8889
$ws.Range("H80").Value = 485.44446
$ws.Range("I80").Value = 452.7143
$ws.Range("J80").Value = 600.0
$ws.Range("K80").Value = 1358.1429
$ws.Range("L80").Value = 1800.0
$ws.Range("M80").Value = -360.1428999999998
$ws.Range("N80").Value = -3796.0
$ws.Range("H83").Value = 485.44446
$ws.Range("I83").Value = 452.7143
$ws.Range("J83").Value = 600.0
$ws.Range("K83").Value = 4074.4287
$ws.Range("L83").Value = 5400.0
$ws.Range("M83").Value = 917.5713000000001
$ws.Range("N83").Value = -15384.0
$ws.Range("H137").Value = 2682.25
$ws.Range("I137").Value = 1645.8334
$ws.Range("K137").Value = 4937.5002
$ws.Range("M137").Value = -2387.5002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 492.83334
$ws.Range("I4").Value = 184.0
$ws.Range("J4").Value = 1110.5
$ws.Range("K4").Value = 184.0
$ws.Range("L4").Value = 1110.5
$ws.Range("M4").Value = -68.0
$ws.Range("N4").Value = -1342.5
$ws.Range("H32").Value = 567.2
$ws.Range("I32").Value = 567.2
$ws.Range("K32").Value = 567.2
$ws.Range("M32").Value = -280.2
$ws.Range("H36").Value = 2600.0
$ws.Range("I36").Value = 2600.0
$ws.Range("K36").Value = 2600.0
$ws.Range("M36").Value = -2254.0
$ws.Range("H61").Value = 6833.1665
$ws.Range("I61").Value = 5749.75
$ws.Range("K61").Value = 5749.75
$ws.Range("M61").Value = -5537.75
$ws.Range("H74").Value = 7017.4
$ws.Range("I74").Value = 7017.4
$ws.Range("K74").Value = 7017.4
$ws.Range("M74").Value = -6143.4
$ws.Range("H77").Value = 7017.4
$ws.Range("I77").Value = 7017.4
$ws.Range("K77").Value = 35087.0
$ws.Range("M77").Value = -30719.0
$ws.Range("H110").Value = 378.54544
$ws.Range("I110").Value = 503.2857
$ws.Range("K110").Value = 503.2857
$ws.Range("M110").Value = 1541.7143
$ws.Range("H132").Value = 1394.3572
$ws.Range("I132").Value = 1394.3572
$ws.Range("K132").Value = 4183.071599999999
$ws.Range("M132").Value = -1653.071599999999
$ws.Range("H136").Value = 6833.1665
$ws.Range("I136").Value = 5749.75
$ws.Range("K136").Value = 17249.25
$ws.Range("M136").Value = -14699.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 7898455.5
$ws.Range("I7").Value = 7861799.5
$ws.Range("K7").Value = 7861799.5
$ws.Range("M7").Value = -7861686.5
$ws.Range("H107").Value = 4370.9443
$ws.Range("I107").Value = 2770.1
$ws.Range("K107").Value = 2770.1
$ws.Range("M107").Value = -850.0999999999999
$ws.Range("H134").Value = 1578.5
$ws.Range("I134").Value = 1309.4445
$ws.Range("J134").Value = 4000.0
$ws.Range("K134").Value = 3928.3335
$ws.Range("L134").Value = 12000.0
$ws.Range("M134").Value = -1393.3335
$ws.Range("N134").Value = -17070.0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2432.5833
$ws.Range("I22").Value = 2327.7144
$ws.Range("K22").Value = 2327.7144
$ws.Range("M22").Value = -1977.7144
$ws.Range("H69").Value = 14595.0
$ws.Range("I69").Value = 7793.3335
$ws.Range("J69").Value = 35000.0
$ws.Range("K69").Value = 7793.3335
$ws.Range("L69").Value = 35000.0
$ws.Range("M69").Value = -7044.3335
$ws.Range("N69").Value = -36498.0
$ws.Range("H72").Value = 14595.0
$ws.Range("I72").Value = 7793.3335
$ws.Range("J72").Value = 35000.0
$ws.Range("K72").Value = 23380.0005
$ws.Range("L72").Value = 105000.0
$ws.Range("M72").Value = -19636.0005
$ws.Range("N72").Value = -112488.0
$ws.Range("H107").Value = 666.1667
$ws.Range("J107").Value = 737.2
$ws.Range("L107").Value = 737.2
$ws.Range("N107").Value = -4577.2
$ws.Range("H122").Value = 764.5625
$ws.Range("I122").Value = 764.5625
$ws.Range("J122").Value = 0.0
$ws.Range("K122").Value = 2293.6875
$ws.Range("L122").Value = 0.0
$ws.Range("M122").Value = 156.3125
$ws.Range("N122").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 101.53333
$ws.Range("J40").Value = 400.66666
$ws.Range("L40").Value = 1602.66664
$ws.Range("N40").Value = -1740.66664
$ws.Range("H86").Value = 612.5
$ws.Range("J86").Value = 750.0
$ws.Range("L86").Value = 2250.0
$ws.Range("N86").Value = -4622.0
$ws.Range("H89").Value = 612.5
$ws.Range("J89").Value = 750.0
$ws.Range("L89").Value = 6750.0
$ws.Range("N89").Value = -18606.0

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 647.2
$ws.Range("I2").Value = 432.6
$ws.Range("J2").Value = 861.8
$ws.Range("K2").Value = 432.6
$ws.Range("L2").Value = 861.8
$ws.Range("M2").Value = -319.6
$ws.Range("N2").Value = -1087.8
$ws.Range("H12").Value = 4166.5
$ws.Range("I12").Value = 2500.0
$ws.Range("K12").Value = 2500.0
$ws.Range("M12").Value = -2360.0
$ws.Range("H70").Value = 10000.0
$ws.Range("J70").Value = 10000.0
$ws.Range("L70").Value = 10000.0
$ws.Range("N70").Value = -10540.0
$ws.Range("H73").Value = 10000.0
$ws.Range("J73").Value = 10000.0
$ws.Range("L73").Value = 10000.0
$ws.Range("N73").Value = -11872.0

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5937.5
$ws.Range("I46").Value = 500.0
$ws.Range("J46").Value = 9200.0
$ws.Range("K46").Value = 500.0
$ws.Range("L46").Value = 9200.0
$ws.Range("M46").Value = -312.0
$ws.Range("N46").Value = -9576.0
$ws.Range("H61").Value = 3093.3333
$ws.Range("I61").Value = 641.0
$ws.Range("J61").Value = 7998.0
$ws.Range("K61").Value = 641.0
$ws.Range("L61").Value = 7998.0
$ws.Range("M61").Value = -439.0
$ws.Range("N61").Value = -8402.0
$ws.Range("H68").Value = 5499.7144
$ws.Range("I68").Value = 2125.0
$ws.Range("K68").Value = 2125.0
$ws.Range("M68").Value = -1376.0
$ws.Range("H71").Value = 5499.7144
$ws.Range("I71").Value = 2125.0
$ws.Range("K71").Value = 10625.0
$ws.Range("M71").Value = -6881.0
$ws.Range("H113").Value = 3093.3333
$ws.Range("I113").Value = 641.0
$ws.Range("J113").Value = 7998.0
$ws.Range("K113").Value = 641.0
$ws.Range("L113").Value = 7998.0
$ws.Range("M113").Value = 1529.0
$ws.Range("N113").Value = -12338.0

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 407.82352
$ws.Range("I107").Value = 370.8125
$ws.Range("K107").Value = 1112.4375
$ws.Range("M107").Value = 807.5625
$ws.Range("H124").Value = 64333.332
$ws.Range("J124").Value = 64333.332
$ws.Range("L124").Value = 64333.332
$ws.Range("N124").Value = -74153.332
